# Update exogenous agent params example values in row 3 (D3:G3)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = 15
$ws.Range("E3").Value = 20
$ws.Range("F3").Value = 10
$ws.Range("G3").Value = 10

# Update the active selection / zoom to match the saved view state
$ws.Range("D4").Select()
$excel.ActiveWindow.Zoom = 235
